$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 237; $r++) {
    $ws.Cells.Item($r, 3).Value = 7312
}

for ($r = 238; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7310
}
